$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data: name, gender, email, hobbies change (city/state unchanged)
$ws.Range("A2").Value = "Mili"
$ws.Range("C2").Value = "Female"
$ws.Range("D2").Value = "mili@gmail.com"
$ws.Range("G2").Value = "Reading,Watching Movies"

# Update the active selection to D2 (matches the saved sheet view state)
$ws.Range("D2").Select()
